# "Generate Report for Handback"
#
# The handback job stamped the localization-status report:
#   - Overview/zh-cn/de-de "Status" columns flip from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - the per-language tables' "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" columns get populated for each localized file
#   - "Latest Target File" becomes a hyperlink (same as the "Source File Name"
#     links already on the sheet)
#   - those newly-populated columns are widened so the content isn't clipped

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (appears on all three sheets: Overview!E:F, and the "Status" column on
#    both the zh-cn and de-de detail sheets)
# ---------------------------------------------------------------------------
foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US") | Out-Null
}

# ---------------------------------------------------------------------------
# 2. Per-language detail sheets (zh-cn, de-de): fill in the handback columns
#    I = Latest Target File   (hyperlink, same display text/target as the
#                               "Source File Name" link in column A)
#    J = Latest Handback File (the generated .xlf file name)
#    K = Latest Handback DateTime
# ---------------------------------------------------------------------------

function Set-HandbackRow(
    $ws,
    [int]$row,
    [string]$mdName,
    [string]$mdUrl,
    [string]$xlfName,
    [string]$handbackDateTime
) {
    $ws.Cells.Item($row, 10).Value = $xlfName               # J: Latest Handback File
    $ws.Cells.Item($row, 11).Value = $handbackDateTime       # K: Latest Handback DateTime
    $ws.Hyperlinks.Add(
        $ws.Cells.Item($row, 9),                             # I: Latest Target File
        $mdUrl,
        "",
        "",
        $mdName
    ) | Out-Null
}

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws2.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws2.Columns.Item(10).ColumnWidth = 39.166666666666664

Set-HandbackRow $ws2 2 `
    "1ee4c6c2-6e30-4c61-980c-0330ef322f42.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c8827fd172b5fa9d7281aca3399ad0ac390ed0b/e2e/1ee4c6c2-6e30-4c61-980c-0330ef322f42.md" `
    "1ee4c6c2-6e30-4c61-980c-0330ef322f42.d775bd04f86a5c9438bee9800cc07797defe2932.zh-cn.xlf" `
    "2016-09-01 23:08:36"

Set-HandbackRow $ws2 3 `
    "1ee4c6c2-6e30-4c61-980c-0330ef322f42.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c8827fd172b5fa9d7281aca3399ad0ac390ed0b/e2e/1ee4c6c2-6e30-4c61-980c-0330ef322f42.md" `
    "1ee4c6c2-6e30-4c61-980c-0330ef322f42.d775bd04f86a5c9438bee9800cc07797defe2932.zh-cn.xlf" `
    "2016-09-01 23:08:36"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws3.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws3.Columns.Item(10).ColumnWidth = 39.166666666666664

Set-HandbackRow $ws3 2 `
    "1ee4c6c2-6e30-4c61-980c-0330ef322f42.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c8827fd172b5fa9d7281aca3399ad0ac390ed0b/e2e/1ee4c6c2-6e30-4c61-980c-0330ef322f42.md" `
    "1ee4c6c2-6e30-4c61-980c-0330ef322f42.d775bd04f86a5c9438bee9800cc07797defe2932.de-de.xlf" `
    "2016-09-01 23:08:44"

Set-HandbackRow $ws3 3 `
    "1ee4c6c2-6e30-4c61-980c-0330ef322f42.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c8827fd172b5fa9d7281aca3399ad0ac390ed0b/e2e/1ee4c6c2-6e30-4c61-980c-0330ef322f42.md" `
    "1ee4c6c2-6e30-4c61-980c-0330ef322f42.d775bd04f86a5c9438bee9800cc07797defe2932.de-de.xlf" `
    "2016-09-01 23:08:44"

# ---------------------------------------------------------------------------
# 3. Overview sheet: widen the "zh-cn" / "de-de" status columns to match
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668
